$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Outline additions under "- Data Abstractions":
#    insert "  - Hexadecimal" before "  - Non-numeric Data", and insert
#    "    - ASCII / Unicode" / "    - Color Models (RGB)" right after it.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "  - Non-numeric Data") {
        $p.Range.InsertParagraphBefore()
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "  - Non-numeric Data") {
        $hexPara = $p.Previous()
        $hexPara.Range.Text = "  - Hexadecimal"
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "  - Non-numeric Data") {
        $p.Range.InsertParagraphAfter()
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "  - Non-numeric Data") {
        $asciiPara = $p.Next()
        $asciiPara.Range.Text = "    - ASCII / Unicode"
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "    - ASCII / Unicode") {
        $p.Range.InsertParagraphAfter()
        break
    }
}
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "    - ASCII / Unicode") {
        $colorPara = $p.Next()
        $colorPara.Range.Text = "    - Color Models (RGB)"
        break
    }
}

# ---------------------------------------------------------------------------
# 2) True/False item: collapse the trailing "." run into the main sentence
#    run (text itself is unchanged).
# ---------------------------------------------------------------------------
$r1 = $d.Content
$null = $r1.Find.Execute(
    "The control unit is one of the three main components of the stored program architecture.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The control unit is one of the three main components of the stored program architecture.",
    2)

# ---------------------------------------------------------------------------
# 3) True/False item: collapse the split "instruction register ... executed."
#    runs down to a single trailing run (text itself is unchanged). Re-typing
#    "instruction register" back onto itself makes the engine coalesce it
#    with all of its same-format neighbour runs to the right. The leading
#    tab (bundled together with "The ") is left untouched so it stays a
#    real <w:tab/> element instead of being flattened into plain text.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$null = $r2.Find.Execute(
    "instruction register",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "instruction register",
    2)

# ---------------------------------------------------------------------------
# 4) True/False item: collapse the " holds the main memory address ..."
#    runs (after "program counter") down to a single run (text unchanged).
#    Scope the touched range tightly to the "micro-program" word so the
#    re-typed text only coalesces with its two immediate same-format
#    neighbour runs and does not reach back across the "program counter"
#    run (which carries a distinguishing w:rsidR and must stay separate).
# ---------------------------------------------------------------------------
$r3 = $d.Content
$null = $r3.Find.Execute("micro-program instruction to be executed", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sub3 = $d.Range($r3.Start, $r3.Start + ("micro-program").Length)
$null = $sub3.Find.Execute(
    "micro-program",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "micro-program",
    2)

# ---------------------------------------------------------------------------
# 5) Mark the first inline picture's run as NoProof (<w:noProof/>).
# ---------------------------------------------------------------------------
if ($d.InlineShapes.Count -ge 1) {
    $shp = $d.InlineShapes(1)
    $shp.Range.NoProofing = 1
}
